# Automatische test-sync: 2025-07-29 21:46:50
# Appends the newest test-mail log entry (Testmail #9) to the "Logs" sheet
# and bumps the matching "Overig" tally on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")

# New row of log data (row 11).
$row = 11
$logs.Cells.Item($row, 1).Value = "Hoi, hebben jullie al iets gehoord?"
$logs.Cells.Item($row, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($row, 3).Value = "Testmail #9: Hoi, hebben jullie al iets gehoord?"
$logs.Cells.Item($row, 4).Value = "Overig"
# Column E (Antwoord) intentionally left blank - no automatic reply for this mail.
$logs.Cells.Item($row, 6).Value = "2025-07-29 21:46:16"
$logs.Cells.Item($row, 7).Value = "Nee"
$logs.Cells.Item($row, 8).Value = "Ja"
$logs.Cells.Item($row, 9).Value = "Nee"
$logs.Cells.Item($row, 10).Value = "Nee"

# Extend the conditional-formatting ranges so the new row is covered too.
$ranges = @("D2:D10", "G2:G10", "H2:H10", "I2:I10", "J2:J10")
$newRanges = @("D2:D11", "G2:G11", "H2:H11", "I2:I11", "J2:J11")
for ($i = 0; $i -lt $ranges.Length; $i++) {
    $fc = $logs.Range($ranges[$i]).FormatConditions
    $target = $logs.Range($newRanges[$i])
    for ($j = 1; $j -le $fc.Count; $j++) {
        $fc.Item($j).ModifyAppliesToRange($target)
    }
}

# Update the Dashboard tally for "Overig" (was 3, now 4).
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Cells.Item(2, 2).Value = 4
